{"js": "// Rewrite the \"Cosmic Symphony\" astronomy essay into the \"Math\" essay,\n// update the author name + email, and append a trailing empty paragraph.\n// Strategy: each paragraph in this document has uniform run formatting\n// (same font/color/size throughout), so we can safely replace the whole\n// paragraph's text in one shot via Range.insertText(..., replace) and the\n// new run Word creates will inherit the formatting that was already there.\n// \"\\u000b\" (vertical tab) is how Word represents a <w:br/> line break in\n// paragraph .text, so it is used for the blank-line dividers inside the\n// long body paragraph.\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nif (paragraphs.items.length < 7) {\n  throw new Error(\"Unexpected document shape: expected at least 7 paragraphs, found \" + paragraphs.items.length);\n}\n\nconst titlePara = paragraphs.items[0];\nconst authorPara = paragraphs.items[1];\nconst emailPara = paragraphs.items[2];\nconst bodyPara = paragraphs.items[4];\nconst summaryBodyPara = paragraphs.items[6];\n\n// 1) Title\ntitlePara.getRange().insertText(\n  \"Math: The Unveiled Language of Patterns and Order\",\n  Word.InsertLocation.replace\n);\n\n// 2) Author name (\"Dr\" + \".\" + \" Olivia Nelson\" -> \"Dr\" + \".\" + \" Amelia Coleman\")\nauthorPara.getRange().insertText(\n  \"Dr. Amelia Coleman\",\n  Word.InsertLocation.replace\n);\n\n// 3) Email address (\"olivianelson@academic\" + \".\" + \"edu\" -> \"ac\" + \".\" + \"integral87@protonmail\" + \".\" + \"ch\")\nemailPara.getRange().insertText(\n  \"ac.integral87@protonmail.ch\",\n  Word.InsertLocation.replace\n);\n\n// 4) Long intro body paragraph (three sentence-groups separated by blank lines made of two <w:br/>)\nconst newBodyText =\n  \"In a world perceived through sensory stimuli, mathematics emerges as a profound realm of thought that penetrates the surface chaos, unveiling the patterns and structure inherent to our universe.\" +\n  \" Beginning with simple arithmetic, which grants us the understanding of quantity and its operations, mathematics expands into a boundless expanse of concepts that explore the forms and changes around us.\" +\n  \" Algebra emerges as a bridge between numbers and geometry, enriching our comprehension of the relationships that govern variables.\" +\n  \" Geometry's axioms and theorems outline the rigid beauty of shapes, angles, and planes, revealing harmonies found in architecture, art, and nature.\" +\n  \" Yet, mathematics does not merely chronicle the known; it delves into the realm of unknown, offering tools to analyze, predict, and solve problems that stretch our intellectual capacities.\" +\n  \"\\u000b\\u000b\" +\n  \"As we venture deeper into this world of numbers, shapes, and relationships, we encounter the intricate relationships that connect different branches of mathematics.\" +\n  \" Analysis unveils the nature of change, unraveling the secrets of calculus and differential equations that drive scientific progress.\" +\n  \" Probability and statistics illuminate the intricate dance of chance, enabling us to make sense of random events and predict future outcomes.\" +\n  \" Amidst this multifaceted tapestry of knowledge, we discover the transformative nature of mathematics as a tool for scientific discovery, innovation, decision-making, and the very act of reasoning itself.\" +\n  \"\\u000b\\u000b\" +\n  \"Just as the luminous night sky inspires awe, so does the realm of mathematics stimulate a sense of wonder, empowering us to decipher the universe's complexities.\" +\n  \" By embracing the challenge and savoring the beauty of mathematical concepts, we navigate the world with greater awareness, understanding, and intellectual fulfillment, recognizing our place within the grand symphony of existence.\";\n\nbodyPara.getRange().insertText(newBodyText, Word.InsertLocation.replace);\n\n// 5) Summary paragraph (heading \"Summary\" is unchanged)\nconst newSummaryText =\n  \"In this essay, we explored the fascinating world of mathematics, delving into its rich history, fundamental concepts, and profound impact on our understanding of the universe.\" +\n  \" From the dawn of arithmetic to the complexities of modern calculus, mathematics serves as a universal language of patterns and order.\" +\n  \" Its branches connect, creating a tapestry of knowledge that drives scientific discovery and technological progress.\" +\n  \" As we unveil the enigmatic tapestry of mathematical concepts, we gain deeper insights into the mysteries of the universe, advancing our intellectual understanding and capacity to make informed decisions.\" +\n  \" Mathematics remains a profound tool, empowering us to navigate the intricacies of a world abundant with patterns and connections, offering new perspectives and unlocking the secrets of our existence.\";\n\nsummaryBodyPara.getRange().insertText(newSummaryText, Word.InsertLocation.replace);\n\nawait context.sync();\n\n// 6) Append a new trailing empty paragraph after the summary paragraph.\nconst paragraphsAfter = context.document.body.paragraphs;\nparagraphsAfter.load(\"items\");\nawait context.sync();\nconst lastPara = paragraphsAfter.items[paragraphsAfter.items.length - 1];\nlastPara.insertParagraph(\"\", Word.InsertLocation.after);\n\nawait context.sync();\n", "ps1": "# Rewrite the \"Cosmic Symphony\" astronomy essay into the \"Math\" essay,\n# update the author name + email, and append a trailing empty paragraph.\n# Each paragraph in this document uses uniform run formatting (same\n# font/color/size throughout the paragraph), so the whole paragraph's\n# text can be replaced in one assignment via Range.Text and the new text\n# inherits the formatting that was already on the paragraph's run(s).\n# [char]11 (vertical tab) is how Word represents a <w:br/> line break in\n# Range.Text, so it is used for the blank-line dividers inside the long\n# body paragraph.\n\n$d = $word.ActiveDocument\n\nif ($d.Paragraphs.Count -lt 7) {\n    throw \"Unexpected document shape: expected at least 7 paragraphs, found $($d.Paragraphs.Count)\"\n}\n\n$nl = [char]11\n\n# 1) Title\n$d.Paragraphs.Item(1).Range.Text = \"Math: The Unveiled Language of Patterns and Order\"\n\n# 2) Author name (\"Dr\" + \".\" + \" Olivia Nelson\" -> \"Dr\" + \".\" + \" Amelia Coleman\")\n$d.Paragraphs.Item(2).Range.Text = \"Dr. Amelia Coleman\"\n\n# 3) Email address (\"olivianelson@academic\" + \".\" + \"edu\" -> \"ac\" + \".\" + \"integral87@protonmail\" + \".\" + \"ch\")\n$d.Paragraphs.Item(3).Range.Text = \"ac.integral87@protonmail.ch\"\n\n# 4) Long intro body paragraph (three sentence-groups separated by blank lines made of two <w:br/>)\n$newBodyText = \"In a world perceived through sensory stimuli, mathematics emerges as a profound realm of thought that penetrates the surface chaos, unveiling the patterns and structure inherent to our universe.\" `\n    + \" Beginning with simple arithmetic, which grants us the understanding of quantity and its operations, mathematics expands into a boundless expanse of concepts that explore the forms and changes around us.\" `\n    + \" Algebra emerges as a bridge between numbers and geometry, enriching our comprehension of the relationships that govern variables.\" `\n    + \" Geometry's axioms and theorems outline the rigid beauty of shapes, angles, and planes, revealing harmonies found in architecture, art, and nature.\" `\n    + \" Yet, mathematics does not merely chronicle the known; it delves into the realm of unknown, offering tools to analyze, predict, and solve problems that stretch our intellectual capacities.\" `\n    + $nl + $nl `\n    + \"As we venture deeper into this world of numbers, shapes, and relationships, we encounter the intricate relationships that connect different branches of mathematics.\" `\n    + \" Analysis unveils the nature of change, unraveling the secrets of calculus and differential equations that drive scientific progress.\" `\n    + \" Probability and statistics illuminate the intricate dance of chance, enabling us to make sense of random events and predict future outcomes.\" `\n    + \" Amidst this multifaceted tapestry of knowledge, we discover the transformative nature of mathematics as a tool for scientific discovery, innovation, decision-making, and the very act of reasoning itself.\" `\n    + $nl + $nl `\n    + \"Just as the luminous night sky inspires awe, so does the realm of mathematics stimulate a sense of wonder, empowering us to decipher the universe's complexities.\" `\n    + \" By embracing the challenge and savoring the beauty of mathematical concepts, we navigate the world with greater awareness, understanding, and intellectual fulfillment, recognizing our place within the grand symphony of existence.\"\n\n$d.Paragraphs.Item(5).Range.Text = $newBodyText\n\n# 5) Summary paragraph (heading \"Summary\" is unchanged)\n$newSummaryText = \"In this essay, we explored the fascinating world of mathematics, delving into its rich history, fundamental concepts, and profound impact on our understanding of the universe.\" `\n    + \" From the dawn of arithmetic to the complexities of modern calculus, mathematics serves as a universal language of patterns and order.\" `\n    + \" Its branches connect, creating a tapestry of knowledge that drives scientific discovery and technological progress.\" `\n    + \" As we unveil the enigmatic tapestry of mathematical concepts, we gain deeper insights into the mysteries of the universe, advancing our intellectual understanding and capacity to make informed decisions.\" `\n    + \" Mathematics remains a profound tool, empowering us to navigate the intricacies of a world abundant with patterns and connections, offering new perspectives and unlocking the secrets of our existence.\"\n\n$d.Paragraphs.Item(7).Range.Text = $newSummaryText\n\n# 6) Append a new trailing empty paragraph after the summary paragraph.\n$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)\n$lastPara.Range.InsertParagraphAfter()\n"}
